$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "53.605.53"
$ws.Range("E2").Value = "  -5.15%  "
$ws.Range("D3").Value = "2.209.63"
$ws.Range("E3").Value = "  -7.24%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.26%  "
$ws.Range("D5").Value = "486.31"
$ws.Range("E5").Value = "  -3.97%  "
$ws.Range("E6").Value = "  -4.44%  "
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  -4.48%  "
$ws.Range("D9").Value = "2.240.05"
$ws.Range("E9").Value = "  -6.54%  "
$ws.Range("E10").Value = "  -6.86%  "
$ws.Range("E11").Value = "  -0.42%  "
$ws.Range("D12").Value = "0.319"
$ws.Range("E12").Value = "  -2.73%  "
$ws.Range("E13").Value = "  -4.17%  "
$ws.Range("D14").Value = "2.630.06"
$ws.Range("E14").Value = "  -6.26%  "
$ws.Range("E15").Value = "  -1.98%  "
$ws.Range("D16").Value = "53.613.64"
$ws.Range("E16").Value = "  -5.08%  "
$ws.Range("E17").Value = "  -4.41%  "
$ws.Range("D18").Value = "2.245.00"
$ws.Range("E18").Value = "  -4.05%  "
$ws.Range("D19").Value = "9.65"
$ws.Range("E19").Value = "  -4.65%  "
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").Value = "296.06"
$ws.Range("E21").Value = "  -4.64%  "
$ws.Range("E22").Value = "  -2.59%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").Value = "63.68"
$ws.Range("E24").Value = "  -4.32%  "
$ws.Range("D25").Value = "'1.00"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  -1.27%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.337.38"
$ws.Range("E27").Value = "  -6.35%  "
$ws.Range("B28").Value = "Kaspa"
$ws.Range("C28").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D28").Value = "0.147"
$ws.Range("E28").Value = "  -0.70%  "
$ws.Range("E29").Value = "  -3.31%  "
$ws.Range("D30").Value = "164.35"
$ws.Range("E30").Value = "  -5.21%  "
$ws.Range("E31").Value = "  -4.59%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'5.80"
$ws.Range("E33").Value = "  -1.07%  "
$ws.Range("B34").Value = "PEPE"
$ws.Range("C34").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D34").Value = "0.0₃0669"
$ws.Range("E34").Value = "  -6.66%  "
$ws.Range("D35").Value = "0.993"
$ws.Range("E35").Value = "  -0.33%  "
$ws.Range("E36").Value = "  -2.32%  "
$ws.Range("E37").Value = "  -2.46%  "
$ws.Range("E38").Value = "  -1.14%  "
$ws.Range("E39").Value = "  +0.97%  "
$ws.Range("E40").Value = "  -4.32%  "
$ws.Range("D41").Value = "35.19"
$ws.Range("E41").Value = "  -3.82%  "
$ws.Range("D42").Value = "0.369"
$ws.Range("E42").Value = "  -1.04%  "
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("D44").Value = "3.29"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("D45").Value = "126.71"
$ws.Range("E45").Value = "  -1.66%  "
$ws.Range("E46").Value = "  -2.95%  "
$ws.Range("E47").Value = "  -2.07%  "
$ws.Range("D48").Value = "0.535"
$ws.Range("E48").Value = "  -5.98%  "
$ws.Range("D49").Value = "232.88"
$ws.Range("E49").Value = "  -3.48%  "
$ws.Range("E50").Value = "  -2.45%  "
$ws.Range("E51").Value = "  -3.55%  "
